$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = "테디노트"

$ws.Range("D6").Value = "[Python - 프로그래머스] 코딩테스트 연습 > 완전탐색 > 피로도"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-%EC%BD%94%EB%94%A9%ED%85%8C%EC%8A%A4%ED%8A%B8-%EC%97%B0%EC%8A%B5-%EC%99%84%EC%A0%84%ED%83%90%EC%83%89-%ED%94%BC%EB%A1%9C%EB%8F%84"

$ws.Range("D32").Value = "[Airflow] Airflow context variable"
$ws.Range("E32").Value = "https://dodonam.tistory.com/401"

$ws.Range("D36").Value = "How advanced is the image semantic segmentation algorithm"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/393"
